$wb = $excel.ActiveWorkbook

# Sheet "展览" (exhibitions) - update 想去人数 (column F) values
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F3").Value = 169
$ws1.Range("F4").Value = 8653
$ws1.Range("F5").Value = 105
$ws1.Range("F11").Value = 178
$ws1.Range("F17").Value = 6123
$ws1.Range("F20").Value = 2253
$ws1.Range("F21").Value = 89
$ws1.Range("F22").Value = 152
$ws1.Range("F23").Value = 244
$ws1.Range("F24").Value = 448

# Sheet "全部类型" (all types) - update 想去人数 (column F) values
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F3").Value = 169
$ws4.Range("F4").Value = 8653
$ws4.Range("F5").Value = 105
$ws4.Range("F13").Value = 178
$ws4.Range("F20").Value = 6124
$ws4.Range("F24").Value = 2253
$ws4.Range("F25").Value = 89
$ws4.Range("F26").Value = 152
$ws4.Range("F27").Value = 244
$ws4.Range("F28").Value = 448
